$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Insert two new columns (Title, Description) after the Question column ---
$ws1.Range("B1:C1").EntireColumn.Insert()

# --- Header row ---
$ws1.Cells.Item(1,1).Value = "Question"
$ws1.Cells.Item(1,2).Value = "Title"
$ws1.Cells.Item(1,3).Value = "Description"
$ws1.Cells.Item(1,4).Value = "Variable_Name"
$ws1.Cells.Item(1,5).Value = "Question_Type"
$ws1.Cells.Item(1,6).Value = "Required"
$ws1.Cells.Item(1,7).Value = "List_Values"
$ws1.Cells.Item(1,8).Value = "If_Condition"
$ws1.Cells.Item(1,9).Value = "Then_Goto"
$ws1.Cells.Item(1,10).Value = "Else_Goto"

# --- Row 2: patient name question ---
$ws1.Cells.Item(2,1).Value = "What is the patient's name?"
$ws1.Cells.Item(2,2).Value = "Patient"
$ws1.Cells.Item(2,3).Value = "Please enter Name…."
$ws1.Cells.Item(2,4).Value = "patient_name"
$ws1.Cells.Item(2,5).Value = "Text"
$ws1.Cells.Item(2,6).Value = $true

# --- Row 3: patient age question ---
$ws1.Cells.Item(3,1).Value = "What is the patient's age?"
$ws1.Cells.Item(3,2).Value = "Age"
$ws1.Cells.Item(3,3).Value = "Please enter Age …."
$ws1.Cells.Item(3,4).Value = "age"
$ws1.Cells.Item(3,5).Value = "Numeric"
$ws1.Cells.Item(3,6).Value = $false

# --- Row 4: patient sex question ---
$ws1.Cells.Item(4,1).Value = "What is the patient's sex?"
$ws1.Cells.Item(4,2).Value = "Sex"
$ws1.Cells.Item(4,3).Value = "Please enter Sex…."
$ws1.Cells.Item(4,4).Value = "sex"
$ws1.Cells.Item(4,5).Value = "Dropdown"
$ws1.Cells.Item(4,6).Value = $false
$ws1.Cells.Item(4,7).Value = "Sheet3"
$ws1.Cells.Item(4,8).Value = "Female"
$ws1.Cells.Item(4,9).Value = "pregnant"
$ws1.Cells.Item(4,10).Value = "onset_date"

# --- Row 5: pregnant question ---
$ws1.Cells.Item(5,1).Value = "Is the patient pregnant?"
$ws1.Cells.Item(5,2).Value = "Pregnant"
$ws1.Cells.Item(5,3).Value = "Please enter Pregnant..."
$ws1.Cells.Item(5,4).Value = "pregnant"
$ws1.Cells.Item(5,5).Value = "Yes/No"
$ws1.Cells.Item(5,6).Value = $false

# --- Row 6: symptom onset question ---
$ws1.Cells.Item(6,1).Value = "When did symptoms start?"
$ws1.Cells.Item(6,2).Value = "symptoms"
$ws1.Cells.Item(6,3).Value = "Please enter symptoms…"
$ws1.Cells.Item(6,4).Value = "onset_date"
$ws1.Cells.Item(6,5).Value = "Date"
$ws1.Cells.Item(6,6).Value = $true

# --- Column widths: the original (bestFit) columns keep their widths after
# the shift caused by the column insert; only the two new columns (Title,
# Description) need an explicit width, matching column A's width.
$ws1.Columns.Item(2).ColumnWidth = 25.276042
$ws1.Columns.Item(3).ColumnWidth = 25.276042

# --- Selection moves to I1 ---
$null = $ws1.Activate()
$null = $ws1.Range("I1").Select()
